$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the column-header shared strings from "<field>_old" / "<field>_new"
#        to "<field>_FV2304" / "<field>_FV2310" (row 1, columns A:U) --------------
$headerMap = @{
    "Segmentname_old"          = "Segmentname_FV2304"
    "Segmentgruppe_old"        = "Segmentgruppe_FV2304"
    "Segment_old"              = "Segment_FV2304"
    "Datenelement_old"         = "Datenelement_FV2304"
    "Segment ID_old"           = "Segment ID_FV2304"
    "Code_old"                 = "Code_FV2304"
    "Qualifier_old"            = "Qualifier_FV2304"
    "Beschreibung_old"         = "Beschreibung_FV2304"
    "Bedingungsausdruck_old"   = "Bedingungsausdruck_FV2304"
    "Bedingung_old"            = "Bedingung_FV2304"
    "Segmentname_new"          = "Segmentname_FV2310"
    "Segmentgruppe_new"        = "Segmentgruppe_FV2310"
    "Segment_new"              = "Segment_FV2310"
    "Datenelement_new"         = "Datenelement_FV2310"
    "Segment ID_new"           = "Segment ID_FV2310"
    "Code_new"                 = "Code_FV2310"
    "Qualifier_new"            = "Qualifier_FV2310"
    "Beschreibung_new"         = "Beschreibung_FV2310"
    "Bedingungsausdruck_new"   = "Bedingungsausdruck_FV2310"
    "Bedingung_new"            = "Bedingung_FV2310"
}

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value2
    if ($headerMap.ContainsKey($current)) {
        $cell.Value = $headerMap[$current]
    }
}

# --- 2. Turn the header row + data into an Excel Table (ListObject) --------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$lo.Name = "Table1"

# --- 3. Freeze the header row (pane split after row 1) ---------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
